$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update CreatedAt timestamp in A1
$ws.Range("A1").Value = "CreatedAt: 2025-11-17T18:06:43"

# Update forecast columns V:Z (hours 20-24) across affected rows
$ws.Cells.Item(4, 22).Value = 145
$ws.Cells.Item(4, 23).Value = 97
$ws.Cells.Item(4, 24).Value = 168.91
$ws.Cells.Item(4, 25).Value = 205.54
$ws.Cells.Item(4, 26).Value = 192.83
$ws.Cells.Item(5, 22).Value = -41.49
$ws.Cells.Item(5, 23).Value = -90.63
$ws.Cells.Item(6, 22).Value = -14.17
$ws.Cells.Item(6, 23).Value = -14.82
$ws.Cells.Item(6, 24).Value = -11.32
$ws.Cells.Item(6, 25).Value = -11.72
$ws.Cells.Item(6, 26).Value = -9.83
$ws.Cells.Item(9, 22).Value = 140.27
$ws.Cells.Item(9, 23).Value = 41.37
$ws.Cells.Item(9, 24).Value = 175.83
$ws.Cells.Item(9, 25).Value = 213.21
$ws.Cells.Item(9, 26).Value = 201.45
$ws.Cells.Item(10, 22).Value = -41.49
$ws.Cells.Item(10, 23).Value = -90.63
$ws.Cells.Item(11, 22).Value = -18.9
$ws.Cells.Item(11, 23).Value = -18.07
$ws.Cells.Item(11, 24).Value = -4.4
$ws.Cells.Item(11, 25).Value = -4.05
$ws.Cells.Item(11, 26).Value = -1.21
$ws.Cells.Item(12, 23).Value = -52.39
$ws.Cells.Item(14, 22).Value = 181.76
$ws.Cells.Item(14, 23).Value = 132
$ws.Cells.Item(14, 24).Value = 212.98
$ws.Cells.Item(14, 25).Value = 213.21
$ws.Cells.Item(14, 26).Value = 212.06
$ws.Cells.Item(15, 23).Value = 0
$ws.Cells.Item(15, 24).Value = 36.98
$ws.Cells.Item(15, 25).Value = 0
$ws.Cells.Item(15, 26).Value = 10.41
$ws.Cells.Item(16, 22).Value = -18.9
$ws.Cells.Item(16, 23).Value = -18.07
$ws.Cells.Item(16, 24).Value = -4.22
$ws.Cells.Item(16, 25).Value = -4.05
$ws.Cells.Item(16, 26).Value = -1.01
$ws.Cells.Item(17, 23).Value = -52.39
$ws.Cells.Item(19, 22).Value = 144
$ws.Cells.Item(19, 23).Value = 168.93
$ws.Cells.Item(19, 24).Value = 167.81
$ws.Cells.Item(19, 25).Value = 204.57
$ws.Cells.Item(19, 26).Value = 191.73
$ws.Cells.Item(20, 22).Value = -40.27
$ws.Cells.Item(20, 23).Value = -16.65
$ws.Cells.Item(21, 22).Value = -16.4
$ws.Cells.Item(21, 23).Value = -16.89
$ws.Cells.Item(21, 24).Value = -12.42
$ws.Cells.Item(21, 25).Value = -12.68
$ws.Cells.Item(21, 26).Value = -10.93
$ws.Cells.Item(24, 22).Value = 142.77
$ws.Cells.Item(24, 23).Value = 94.94
$ws.Cells.Item(24, 24).Value = 167.81
$ws.Cells.Item(24, 25).Value = 204.57
$ws.Cells.Item(24, 26).Value = 191.73
$ws.Cells.Item(25, 22).Value = -41.49
$ws.Cells.Item(25, 23).Value = -90.63
$ws.Cells.Item(26, 22).Value = -16.4
$ws.Cells.Item(26, 23).Value = -16.89
$ws.Cells.Item(26, 24).Value = -12.42
$ws.Cells.Item(26, 25).Value = -12.68
$ws.Cells.Item(26, 26).Value = -10.93
$ws.Cells.Item(29, 22).Value = 141.99
$ws.Cells.Item(29, 23).Value = 166.74
$ws.Cells.Item(29, 24).Value = 166.41
$ws.Cells.Item(29, 25).Value = 202.86
$ws.Cells.Item(29, 26).Value = 190.29
$ws.Cells.Item(30, 22).Value = -40.27
$ws.Cells.Item(30, 23).Value = -16.65
$ws.Cells.Item(31, 22).Value = -18.41
$ws.Cells.Item(31, 23).Value = -19.07
$ws.Cells.Item(31, 24).Value = -13.81
$ws.Cells.Item(31, 25).Value = -14.4
$ws.Cells.Item(31, 26).Value = -12.37
$ws.Cells.Item(34, 22).Value = 170
$ws.Cells.Item(34, 23).Value = 157.59
$ws.Cells.Item(34, 24).Value = 218.29
$ws.Cells.Item(34, 25).Value = 219.23
$ws.Cells.Item(34, 26).Value = 219.12
$ws.Cells.Item(35, 23).Value = 0
$ws.Cells.Item(35, 24).Value = 36.98
$ws.Cells.Item(35, 25).Value = 0
$ws.Cells.Item(35, 26).Value = 10.41
$ws.Cells.Item(36, 22).Value = -19.56
$ws.Cells.Item(36, 23).Value = -17.4
$ws.Cells.Item(36, 24).Value = 1.09
$ws.Cells.Item(36, 25).Value = 1.97
$ws.Cells.Item(36, 26).Value = 6.05
$ws.Cells.Item(37, 22).Value = -11.11
$ws.Cells.Item(37, 23).Value = -27.47
$ws.Cells.Item(39, 22).Value = 145
$ws.Cells.Item(39, 23).Value = 97
$ws.Cells.Item(39, 24).Value = 168.91
$ws.Cells.Item(39, 25).Value = 205.54
$ws.Cells.Item(39, 26).Value = 192.83
$ws.Cells.Item(40, 22).Value = -41.49
$ws.Cells.Item(40, 23).Value = -90.63
$ws.Cells.Item(41, 22).Value = -14.17
$ws.Cells.Item(41, 23).Value = -14.82
$ws.Cells.Item(41, 24).Value = -11.32
$ws.Cells.Item(41, 25).Value = -11.72
$ws.Cells.Item(41, 26).Value = -9.83
$ws.Cells.Item(44, 22).Value = 199.07
$ws.Cells.Item(44, 23).Value = 199.66
$ws.Cells.Item(44, 24).Value = 177.39
$ws.Cells.Item(44, 25).Value = 214.47
$ws.Cells.Item(44, 26).Value = 201.65
$ws.Cells.Item(46, 22).Value = -1.59
$ws.Cells.Item(46, 23).Value = -2.8
$ws.Cells.Item(46, 24).Value = -2.84
$ws.Cells.Item(46, 26).Value = -1.01
$ws.Cells.Item(49, 22).Value = 201.07
$ws.Cells.Item(49, 23).Value = 214.01
$ws.Cells.Item(49, 24).Value = 189.71
$ws.Cells.Item(49, 25).Value = 226.78
$ws.Cells.Item(49, 26).Value = 209.14
$ws.Cells.Item(51, 22).Value = 0.4
$ws.Cells.Item(51, 23).Value = 11.56
$ws.Cells.Item(51, 24).Value = 9.49
$ws.Cells.Item(51, 25).Value = 9.52
$ws.Cells.Item(51, 26).Value = 6.48
$ws.Cells.Item(54, 22).Value = 192.58
$ws.Cells.Item(54, 23).Value = 201.45
$ws.Cells.Item(54, 24).Value = 181.86
$ws.Cells.Item(54, 25).Value = 221.47
$ws.Cells.Item(54, 26).Value = 211.1
$ws.Cells.Item(56, 22).Value = -8.09
$ws.Cells.Item(56, 23).Value = -1.01
$ws.Cells.Item(56, 24).Value = 1.64
$ws.Cells.Item(56, 25).Value = 4.21
$ws.Cells.Item(56, 26).Value = 8.44
$ws.Cells.Item(59, 22).Value = 207.09
$ws.Cells.Item(59, 23).Value = 208.08
$ws.Cells.Item(59, 24).Value = 184.28
$ws.Cells.Item(59, 25).Value = 222.6
$ws.Cells.Item(59, 26).Value = 210.01
$ws.Cells.Item(61, 22).Value = 6.42
$ws.Cells.Item(61, 24).Value = 4.05
$ws.Cells.Item(61, 25).Value = 5.34
$ws.Cells.Item(61, 26).Value = 7.35
$ws.Cells.Item(64, 22).Value = 211.23
$ws.Cells.Item(64, 23).Value = 212
$ws.Cells.Item(64, 24).Value = 187.34
$ws.Cells.Item(64, 25).Value = 226.07
$ws.Cells.Item(64, 26).Value = 213.55
$ws.Cells.Item(66, 22).Value = 10.56
$ws.Cells.Item(66, 23).Value = 9.539999999999999
$ws.Cells.Item(66, 24).Value = 7.12
$ws.Cells.Item(66, 25).Value = 8.82
$ws.Cells.Item(66, 26).Value = 10.89
$ws.Cells.Item(69, 22).Value = 212.35
$ws.Cells.Item(69, 23).Value = 212.89
$ws.Cells.Item(69, 24).Value = 187.93
$ws.Cells.Item(69, 25).Value = 227.02
$ws.Cells.Item(69, 26).Value = 215.37
$ws.Cells.Item(71, 22).Value = 11.68
$ws.Cells.Item(71, 23).Value = 10.43
$ws.Cells.Item(71, 24).Value = 7.71
$ws.Cells.Item(71, 25).Value = 9.76
$ws.Cells.Item(71, 26).Value = 12.71
$ws.Cells.Item(74, 22).Value = 208.38
$ws.Cells.Item(74, 23).Value = 210.46
$ws.Cells.Item(74, 24).Value = 186.18
$ws.Cells.Item(74, 25).Value = 224.67
$ws.Cells.Item(74, 26).Value = 211.54
$ws.Cells.Item(76, 22).Value = 7.71
$ws.Cells.Item(76, 23).Value = 8
$ws.Cells.Item(76, 24).Value = 5.96
$ws.Cells.Item(76, 25).Value = 7.41
$ws.Cells.Item(76, 26).Value = 8.880000000000001
$ws.Cells.Item(79, 22).Value = 209.33
$ws.Cells.Item(79, 23).Value = 211.04
$ws.Cells.Item(79, 24).Value = 187.09
$ws.Cells.Item(79, 25).Value = 225.81
$ws.Cells.Item(79, 26).Value = 212.65
$ws.Cells.Item(81, 22).Value = 8.67
$ws.Cells.Item(81, 23).Value = 8.59
$ws.Cells.Item(81, 24).Value = 6.86
$ws.Cells.Item(81, 25).Value = 8.550000000000001
$ws.Cells.Item(81, 26).Value = 9.99
$ws.Cells.Item(84, 22).Value = 186.32
$ws.Cells.Item(84, 23).Value = 193.74
$ws.Cells.Item(84, 24).Value = 174.3
$ws.Cells.Item(84, 25).Value = 211.55
$ws.Cells.Item(84, 26).Value = 205.96
$ws.Cells.Item(86, 22).Value = -14.35
$ws.Cells.Item(86, 23).Value = -8.720000000000001
$ws.Cells.Item(86, 24).Value = -5.93
$ws.Cells.Item(86, 25).Value = -5.71
$ws.Cells.Item(86, 26).Value = 3.3
$ws.Cells.Item(89, 22).Value = 140.77
$ws.Cells.Item(89, 23).Value = 92.75
$ws.Cells.Item(89, 24).Value = 166.41
$ws.Cells.Item(89, 25).Value = 202.86
$ws.Cells.Item(89, 26).Value = 190.29
$ws.Cells.Item(90, 22).Value = -41.49
$ws.Cells.Item(90, 23).Value = -90.63
$ws.Cells.Item(91, 22).Value = -18.41
$ws.Cells.Item(91, 23).Value = -19.07
$ws.Cells.Item(91, 24).Value = -13.81
$ws.Cells.Item(91, 25).Value = -14.4
$ws.Cells.Item(91, 26).Value = -12.37
